$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the surviving data rows (2-6) with their new values.
$ws.Range("A2").Value = "COMBUSTION_FIJA"
$ws.Range("B2").Value = "GAS_NATURAL"
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = "MENSUAL"
$ws.Range("E2").Value = 44652

$ws.Range("A3").Value = "LOGISTICA_PRODUCTOS_RESIDUOS"
$ws.Range("B3").Value = "CATEGORIA"
$ws.Range("C3").Value = "MATERIA_PRIMA"
$ws.Range("D3").Value = "MENSUAL"
$ws.Range("E3").Value = 44682

$ws.Range("A4").Value = "LOGISTICA_PRODUCTOS_RESIDUOS"
$ws.Range("B4").Value = "MEDIO_TRANSPORTE"
$ws.Range("C4").Value = "CAMION_CARGA"
$ws.Range("D4").Value = "MENSUAL"
$ws.Range("E4").Value = 44682

$ws.Range("A5").Value = "LOGISTICA_PRODUCTOS_RESIDUOS"
$ws.Range("B5").Value = "DISTANCIA"
$ws.Range("C5").Value = 80
$ws.Range("D5").Value = "MENSUAL"
$ws.Range("E5").Value = 44682

$ws.Range("A6").Value = "LOGISTICA_PRODUCTOS_RESIDUOS"
$ws.Range("B6").Value = "PESO"
$ws.Range("C6").Value = 800
$ws.Range("D6").Value = "MENSUAL"
$ws.Range("E6").Value = 44682

# Rows 7-11 no longer exist in the new version of the sheet - delete them
# entirely so the used range shrinks back down to A1:E6.
$ws.Range("A7:E11").Delete()

# Remove the bespoke "best fit" column widths the old layout had.
$ws.Columns("A:E").AutoFit()
$ws.Columns("A:E").ColumnWidth = 8.43

# Restore the default selection/view state that Excel writes after a
# fresh edit (was G12, now H2).
$ws.Range("H2").Select()
